# Generate Report for Archive
#
# The localization status for this file moved on from "Ready for handoff"
# to "In Translation" - update the status cells that held the old text on
# every sheet (Overview + each locale table), then let the status columns
# re-shrink to fit the new (shorter) text, same as Excel would do when the
# column was already sized to its contents.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target OOXML column "width" for the shrunk status columns is ~13.41
# character-units. Excel's ColumnWidth property only accepts values on its
# MDW pixel grid, so 12.5 is the nearest settable value that still rounds
# down to the same stored width as the narrower "In Translation" text.
$statusColWidth = 12.5

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
